$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 701
$ws.Range("I9").Value = 949.2778
$ws.Range("J9").Value = 62.57143
$ws.Range("K9").Value = 949.2778
$ws.Range("L9").Value = 62.57143
$ws.Range("M9").Value = -780.2778
$ws.Range("N9").Value = -400.57143
$ws.Range("H19").Value = 5627.095
$ws.Range("I19").Value = 10326.5
$ws.Range("J19").Value = 1354.909
$ws.Range("K19").Value = 10326.5
$ws.Range("L19").Value = 1354.909
$ws.Range("M19").Value = -10151.5
$ws.Range("N19").Value = -1704.909
$ws.Range("H33").Value = 1018.1818
$ws.Range("J33").Value = 472
$ws.Range("L33").Value = 472
$ws.Range("N33").Value = -930
$ws.Range("H76").Value = 9372.137000000001
$ws.Range("I76").Value = 15309.667
$ws.Range("J76").Value = 5261.5386
$ws.Range("K76").Value = 15309.667
$ws.Range("L76").Value = 5261.5386
$ws.Range("M76").Value = -14994.667
$ws.Range("N76").Value = -5891.5386
$ws.Range("H79").Value = 9372.137000000001
$ws.Range("I79").Value = 15309.667
$ws.Range("J79").Value = 5261.5386
$ws.Range("K79").Value = 15309.667
$ws.Range("L79").Value = 5261.5386
$ws.Range("M79").Value = -14217.667
$ws.Range("N79").Value = -7445.5386
$ws.Range("H113").Value = 3650.923
$ws.Range("I113").Value = 3106.25
$ws.Range("K113").Value = 3106.25
$ws.Range("M113").Value = 147.75
$ws.Range("H132").Value = 2771.9473
$ws.Range("I132").Value = 1808.9778
$ws.Range("J132").Value = 6383.0835
$ws.Range("K132").Value = 5426.9334
$ws.Range("L132").Value = 19149.2505
$ws.Range("M132").Value = -2896.9334
$ws.Range("N132").Value = -24209.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1772.2745
$ws.Range("I74").Value = 1138.9269
$ws.Range("J74").Value = 4369
$ws.Range("K74").Value = 1138.9269
$ws.Range("L74").Value = 4369
$ws.Range("M74").Value = -264.9268999999999
$ws.Range("N74").Value = -6117
$ws.Range("H77").Value = 1772.2745
$ws.Range("I77").Value = 1138.9269
$ws.Range("J77").Value = 4369
$ws.Range("K77").Value = 5694.6345
$ws.Range("L77").Value = 21845
$ws.Range("M77").Value = -1326.6345
$ws.Range("N77").Value = -30581
$ws.Range("H122").Value = 3160.2368
$ws.Range("I122").Value = 3101.125
$ws.Range("K122").Value = 9303.375
$ws.Range("M122").Value = -6853.375
$ws.Range("H132").Value = 27297.453
$ws.Range("I132").Value = 51330.383
$ws.Range("K132").Value = 153991.149
$ws.Range("M132").Value = -151461.149

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1106.1923
$ws.Range("I94").Value = 824.5789
$ws.Range("J94").Value = 1870.5714
$ws.Range("K94").Value = 824.5789
$ws.Range("L94").Value = 1870.5714
$ws.Range("M94").Value = -373.5789
$ws.Range("N94").Value = -2772.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 23785.5
$ws.Range("J47").Value = 23785.5
$ws.Range("L47").Value = 23785.5
$ws.Range("N47").Value = -24917.5
$ws.Range("H58").Value = 1821.5518
$ws.Range("I58").Value = 1034
$ws.Range("J58").Value = 2556.6
$ws.Range("K58").Value = 1034
$ws.Range("L58").Value = 2556.6
$ws.Range("M58").Value = -831
$ws.Range("N58").Value = -2962.6
$ws.Range("H62").Value = 1919744.4
$ws.Range("J62").Value = 4367.136
$ws.Range("L62").Value = 4367.136
$ws.Range("N62").Value = -5615.136
$ws.Range("H65").Value = 1919744.4
$ws.Range("J65").Value = 4367.136
$ws.Range("L65").Value = 21835.68
$ws.Range("N65").Value = -28075.68
$ws.Range("H132").Value = 1695.5781
$ws.Range("I132").Value = 974.53656
$ws.Range("J132").Value = 2980.913
$ws.Range("K132").Value = 2923.60968
$ws.Range("L132").Value = 8942.739
$ws.Range("M132").Value = -393.60968
$ws.Range("N132").Value = -14002.739
$ws.Range("H134").Value = 1360.0834
$ws.Range("I134").Value = 871
$ws.Range("J134").Value = 2547.8572
$ws.Range("K134").Value = 2613
$ws.Range("L134").Value = 7643.571599999999
$ws.Range("M134").Value = -78
$ws.Range("N134").Value = -12713.5716
$ws.Range("H136").Value = 1821.5518
$ws.Range("I136").Value = 1034
$ws.Range("J136").Value = 2556.6
$ws.Range("K136").Value = 3102
$ws.Range("L136").Value = 7669.799999999999
$ws.Range("M136").Value = -552
$ws.Range("N136").Value = -12769.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 313.8125
$ws.Range("I68").Value = 262.57144
$ws.Range("J68").Value = 353.66666
$ws.Range("K68").Value = 787.71432
$ws.Range("L68").Value = 1060.99998
$ws.Range("M68").Value = 23.28567999999996
$ws.Range("N68").Value = -2682.99998
$ws.Range("H71").Value = 313.8125
$ws.Range("I71").Value = 262.57144
$ws.Range("J71").Value = 353.66666
$ws.Range("K71").Value = 2363.14296
$ws.Range("L71").Value = 3182.99994
$ws.Range("M71").Value = 1692.85704
$ws.Range("N71").Value = -11294.99994
$ws.Range("H131").Value = 1555.5116
$ws.Range("I131").Value = 441
$ws.Range("J131").Value = 1850.5294
$ws.Range("K131").Value = 1323
$ws.Range("L131").Value = 5551.5882
$ws.Range("M131").Value = 3717
$ws.Range("N131").Value = -15631.5882
$ws.Range("H139").Value = 2098.5334
$ws.Range("J139").Value = 2948.1667
$ws.Range("L139").Value = 8844.500100000001
$ws.Range("N139").Value = -19124.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6472
$ws.Range("I80").Value = 6617.0586
$ws.Range("J80").Value = 4006
$ws.Range("K80").Value = 6617.0586
$ws.Range("L80").Value = 4006
$ws.Range("M80").Value = -5619.0586
$ws.Range("N80").Value = -6002
$ws.Range("H83").Value = 6472
$ws.Range("I83").Value = 6617.0586
$ws.Range("J83").Value = 4006
$ws.Range("K83").Value = 33085.29300000001
$ws.Range("L83").Value = 20030
$ws.Range("M83").Value = -28093.29300000001
$ws.Range("N83").Value = -30014
$ws.Range("H122").Value = 1849.7858
$ws.Range("I122").Value = 1971
$ws.Range("J122").Value = 1728.5714
$ws.Range("K122").Value = 5913
$ws.Range("L122").Value = 5185.7142
$ws.Range("M122").Value = -3463
$ws.Range("N122").Value = -10085.7142
$ws.Range("H126").Value = 3411.6667
$ws.Range("I126").Value = 3355.2632
$ws.Range("J126").Value = 3474.7058
$ws.Range("K126").Value = 10065.7896
$ws.Range("L126").Value = 10424.1174
$ws.Range("M126").Value = -7595.7896
$ws.Range("N126").Value = -15364.1174
$ws.Range("H127").Value = 21500
$ws.Range("J127").Value = 21500
$ws.Range("L127").Value = 21500
$ws.Range("N127").Value = -31420
$ws.Range("H132").Value = 3959.3416
$ws.Range("I132").Value = 3839.484
$ws.Range("J132").Value = 4330.9
$ws.Range("K132").Value = 11518.452
$ws.Range("L132").Value = 12992.7
$ws.Range("M132").Value = -8988.451999999999
$ws.Range("N132").Value = -18052.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3465.3076
$ws.Range("I7").Value = 3960
$ws.Range("J7").Value = 3156.125
$ws.Range("K7").Value = 3960
$ws.Range("L7").Value = 3156.125
$ws.Range("M7").Value = -3848
$ws.Range("N7").Value = -3380.125
$ws.Range("H46").Value = 1598.5385
$ws.Range("I46").Value = 1838.1
$ws.Range("J46").Value = 800
$ws.Range("K46").Value = 1838.1
$ws.Range("L46").Value = 800
$ws.Range("M46").Value = -1650.1
$ws.Range("N46").Value = -1176
$ws.Range("H53").Value = 11637.25
$ws.Range("I53").Value = 7049
$ws.Range("J53").Value = 13166.667
$ws.Range("K53").Value = 7049
$ws.Range("L53").Value = 13166.667
$ws.Range("M53").Value = -6531
$ws.Range("N53").Value = -14202.667
$ws.Range("H126").Value = 3465.3076
$ws.Range("I126").Value = 3960
$ws.Range("J126").Value = 3156.125
$ws.Range("K126").Value = 11880
$ws.Range("L126").Value = 9468.375
$ws.Range("M126").Value = -9410
$ws.Range("N126").Value = -14408.375
$ws.Range("H132").Value = 6184.75
$ws.Range("I132").Value = 2211.5908
$ws.Range("J132").Value = 9546.654
$ws.Range("K132").Value = 6634.7724
$ws.Range("L132").Value = 28639.962
$ws.Range("M132").Value = -4104.7724
$ws.Range("N132").Value = -33699.962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 1005000
$ws.Range("J53").Value = 10000
$ws.Range("L53").Value = 10000
$ws.Range("N53").Value = -11214
$ws.Range("H122").Value = 27779396
$ws.Range("I122").Value = 40001680
$ws.Range("J122").Value = 1475.909
$ws.Range("K122").Value = 120005040
$ws.Range("L122").Value = 4427.727000000001
$ws.Range("M122").Value = -120002590
$ws.Range("N122").Value = -9327.727000000001
$ws.Range("H132").Value = 1823.775
$ws.Range("I132").Value = 1323.3182
$ws.Range("J132").Value = 2435.4443
$ws.Range("K132").Value = 3969.9546
$ws.Range("L132").Value = 7306.3329
$ws.Range("M132").Value = -1439.9546
$ws.Range("N132").Value = -12366.3329
